$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1058.3704
$ws.Range("J112").Value = 1068.3077
$ws.Range("L112").Value = 3204.9231
$ws.Range("N112").Value = -5420.9231
$ws.Range("H116").Value = 4434
$ws.Range("J116").Value = 4715.143
$ws.Range("L116").Value = 4715.143
$ws.Range("N116").Value = -11599.143
$ws.Range("H129").Value = 176434.44
$ws.Range("J129").Value = 186201.17
$ws.Range("L129").Value = 558603.51
$ws.Range("N129").Value = -568603.51
$ws.Range("H132").Value = 2326.2554
$ws.Range("I132").Value = 2548.2
$ws.Range("K132").Value = 7644.599999999999
$ws.Range("M132").Value = -5114.599999999999
$ws.Range("H136").Value = 47519.332
$ws.Range("J136").Value = 47519.332
$ws.Range("L136").Value = 47519.332
$ws.Range("N136").Value = -57719.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 799.2727
$ws.Range("I2").Value = 816.8946999999999
$ws.Range("J2").Value = 687.6667
$ws.Range("K2").Value = 816.8946999999999
$ws.Range("L2").Value = 687.6667
$ws.Range("M2").Value = -703.8946999999999
$ws.Range("N2").Value = -913.6667
$ws.Range("H44").Value = 29016.334
$ws.Range("J44").Value = 29016.334
$ws.Range("L44").Value = 29016.334
$ws.Range("N44").Value = -29992.334
$ws.Range("H45").Value = 3161.1292
$ws.Range("I45").Value = 3068.6428
$ws.Range("J45").Value = 3237.2942
$ws.Range("K45").Value = 3068.6428
$ws.Range("L45").Value = 3237.2942
$ws.Range("M45").Value = -2691.6428
$ws.Range("N45").Value = -3991.2942
$ws.Range("H55").Value = 28788.25
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 28788.25
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 28788.25
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -29418.25
$ws.Range("H116").Value = 799.2727
$ws.Range("I116").Value = 816.8946999999999
$ws.Range("J116").Value = 687.6667
$ws.Range("K116").Value = 816.8946999999999
$ws.Range("L116").Value = 687.6667
$ws.Range("M116").Value = 1477.1053
$ws.Range("N116").Value = -5275.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 799.2727
$ws.Range("I3").Value = 816.8946999999999
$ws.Range("J3").Value = 687.6667
$ws.Range("K3").Value = 816.8946999999999
$ws.Range("L3").Value = 687.6667
$ws.Range("M3").Value = -702.8946999999999
$ws.Range("N3").Value = -915.6667
$ws.Range("H20").Value = 2729.3
$ws.Range("I20").Value = 3153
$ws.Range("J20").Value = 1942.4286
$ws.Range("K20").Value = 3153
$ws.Range("L20").Value = 1942.4286
$ws.Range("M20").Value = -2906
$ws.Range("N20").Value = -2436.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 18019.934
$ws.Range("I58").Value = 1340.4584
$ws.Range("K58").Value = 1340.4584
$ws.Range("M58").Value = -1137.4584
$ws.Range("H132").Value = 2104.7273
$ws.Range("I132").Value = 1580.8485
$ws.Range("J132").Value = 3676.3635
$ws.Range("K132").Value = 4742.5455
$ws.Range("L132").Value = 11029.0905
$ws.Range("M132").Value = -2212.5455
$ws.Range("N132").Value = -16089.0905
$ws.Range("H134").Value = 846.7234
$ws.Range("I134").Value = 786.1111
$ws.Range("J134").Value = 1045.091
$ws.Range("K134").Value = 2358.3333
$ws.Range("L134").Value = 3135.273
$ws.Range("M134").Value = 176.6667000000002
$ws.Range("N134").Value = -8205.272999999999
$ws.Range("H136").Value = 18019.934
$ws.Range("I136").Value = 1340.4584
$ws.Range("K136").Value = 4021.3752
$ws.Range("M136").Value = -1471.3752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1878.6
$ws.Range("I98").Value = 1331
$ws.Range("J98").Value = 2700
$ws.Range("K98").Value = 3993
$ws.Range("L98").Value = 8100
$ws.Range("M98").Value = -2495
$ws.Range("N98").Value = -11096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20813.414
$ws.Range("I132").Value = 3690.04
$ws.Range("K132").Value = 11070.12
$ws.Range("M132").Value = -8540.119999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1691.0416
$ws.Range("I46").Value = 1722.0454
$ws.Range("J46").Value = 1350
$ws.Range("K46").Value = 1722.0454
$ws.Range("L46").Value = 1350
$ws.Range("M46").Value = -1534.0454
$ws.Range("N46").Value = -1726
$ws.Range("H55").Value = 891.3333
$ws.Range("I55").Value = 1600
$ws.Range("J55").Value = 182.66667
$ws.Range("K55").Value = 1600
$ws.Range("L55").Value = 182.66667
$ws.Range("M55").Value = -1427
$ws.Range("N55").Value = -528.6666700000001
$ws.Range("H68").Value = 2514.1428
$ws.Range("I68").Value = 1902
$ws.Range("J68").Value = 2616.1667
$ws.Range("K68").Value = 1902
$ws.Range("L68").Value = 2616.1667
$ws.Range("M68").Value = -1153
$ws.Range("N68").Value = -4114.1667
$ws.Range("H71").Value = 2514.1428
$ws.Range("I71").Value = 1902
$ws.Range("J71").Value = 2616.1667
$ws.Range("K71").Value = 9510
$ws.Range("L71").Value = 13080.8335
$ws.Range("M71").Value = -5766
$ws.Range("N71").Value = -20568.8335
$ws.Range("H132").Value = 525759.6
$ws.Range("I132").Value = 804366.25
$ws.Range("J132").Value = 3372.125
$ws.Range("K132").Value = 2413098.75
$ws.Range("L132").Value = 10116.375
$ws.Range("M132").Value = -2410568.75
$ws.Range("N132").Value = -15176.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4182.091
$ws.Range("I62").Value = 3100
$ws.Range("J62").Value = 5083.8335
$ws.Range("K62").Value = 3100
$ws.Range("L62").Value = 5083.8335
$ws.Range("M62").Value = -2476
$ws.Range("N62").Value = -6331.8335
$ws.Range("H65").Value = 4182.091
$ws.Range("I65").Value = 3100
$ws.Range("J65").Value = 5083.8335
$ws.Range("K65").Value = 15500
$ws.Range("L65").Value = 25419.1675
$ws.Range("M65").Value = -12380
$ws.Range("N65").Value = -31659.1675
$ws.Range("H81").Value = 1693.8572
$ws.Range("I81").Value = 1067.2727
$ws.Range("K81").Value = 2134.5454
$ws.Range("M81").Value = -1073.5454
$ws.Range("H84").Value = 1693.8572
$ws.Range("I84").Value = 1067.2727
$ws.Range("K84").Value = 10672.727
$ws.Range("M84").Value = -5368.726999999999
$ws.Range("H100").Value = 466
$ws.Range("I100").Value = 512
$ws.Range("J100").Value = 351
$ws.Range("K100").Value = 1024
$ws.Range("L100").Value = 702
$ws.Range("M100").Value = -483
$ws.Range("N100").Value = -1784
$ws.Range("H135").Value = 47739.11
$ws.Range("J135").Value = 47739.11
$ws.Range("L135").Value = 47739.11
$ws.Range("N135").Value = -57879.11
$ws.Range("H136").Value = 17859312
$ws.Range("I136").Value = 24391136
$ws.Range("J136").Value = 5659.2
$ws.Range("K136").Value = 73173408
$ws.Range("L136").Value = 16977.6
$ws.Range("M136").Value = -73170858
$ws.Range("N136").Value = -22077.6
$ws.Range("H141").Value = 62969.168
$ws.Range("J141").Value = 62969.168
$ws.Range("L141").Value = 62969.168
$ws.Range("N141").Value = -73329.16800000001
